$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6122626666666666
$ws.Range("H2").Value = 1.836788
$ws.Range("I2").Value = 0.006779070576782467
$ws.Range("J2").Value = 0.006779070576782467
$ws.Range("M2").Value = 41.60833666666667
$ws.Range("N2").Value = 124.82501
$ws.Range("O2").Value = 0.5886423873735626
$ws.Range("P2").Value = 0.5886423873735626
$ws.Range("Q2").Value = 25.47523116309777
$ws.Range("R2").Value = 229.27708046788
$ws.Range("S2").Value = 0.003990448288491105
$ws.Range("T2").Value = 0.003990448288491105
$ws.Range("G3").Value = 0.6122626666666666
$ws.Range("H3").Value = 1.836788
$ws.Range("I3").Value = 0.006779070576782467
$ws.Range("J3").Value = 0.006779070576782467
$ws.Range("O3").Value = 0.2506218293658061
$ws.Range("P3").Value = 0.2506218293658061
$ws.Range("Q3").Value = 10.84639702230711
$ws.Range("R3").Value = 97.617573200764
$ws.Range("S3").Value = 0.001698983069353132
$ws.Range("T3").Value = 0.001698983069353132
$ws.Range("G4").Value = 0.6122626666666666
$ws.Range("H4").Value = 1.836788
$ws.Range("I4").Value = 0.006779070576782467
$ws.Range("J4").Value = 0.006779070576782467
$ws.Range("M4").Value = 6.935318333333334
$ws.Range("N4").Value = 20.805955
$ws.Range("O4").Value = 0.09811549001908282
$ws.Range("P4").Value = 0.09811549001908279
$ws.Range("Q4").Value = 4.246236496948889
$ws.Range("R4").Value = 38.21612847254
$ws.Range("S4").Value = 0.0006651318315149581
$ws.Range("T4").Value = 0.0006651318315149579
$ws.Range("G5").Value = 0.6122626666666666
$ws.Range("H5").Value = 1.836788
$ws.Range("I5").Value = 0.006779070576782467
$ws.Range("J5").Value = 0.006779070576782467
$ws.Range("M5").Value = 4.426331333333334
$ws.Range("N5").Value = 13.278994
$ws.Range("O5").Value = 0.06262029324154843
$ws.Range("P5").Value = 0.06262029324154841
$ws.Range("Q5").Value = 2.710077425696889
$ws.Range("R5").Value = 24.390696831272
$ws.Range("S5").Value = 0.0004245073874232709
$ws.Range("T5").Value = 0.0004245073874232708
$ws.Range("I6").Value = 0.003538518590750013
$ws.Range("J6").Value = 0.003538518590750013
$ws.Range("M6").Value = 41.60833666666667
$ws.Range("N6").Value = 124.82501
$ws.Range("O6").Value = 0.5886423873735626
$ws.Range("P6").Value = 0.5886423873735626
$ws.Range("Q6").Value = 13.29748349029
$ws.Range("R6").Value = 119.67735141261
$ws.Range("S6").Value = 0.002082922031024822
$ws.Range("T6").Value = 0.002082922031024822
$ws.Range("I7").Value = 0.003538518590750013
$ws.Range("J7").Value = 0.003538518590750013
$ws.Range("O7").Value = 0.2506218293658061
$ws.Range("P7").Value = 0.2506218293658061
$ws.Range("S7").Value = 0.0008868300024586826
$ws.Range("T7").Value = 0.0008868300024586826
$ws.Range("I8").Value = 0.003538518590750013
$ws.Range("J8").Value = 0.003538518590750013
$ws.Range("M8").Value = 6.935318333333334
$ws.Range("N8").Value = 20.805955
$ws.Range("O8").Value = 0.09811549001908282
$ws.Range("P8").Value = 0.09811549001908279
$ws.Range("Q8").Value = 2.216437580195
$ws.Range("R8").Value = 19.947938221755
$ws.Range("S8").Value = 0.0003471834854730719
$ws.Range("T8").Value = 0.0003471834854730718
$ws.Range("I9").Value = 0.003538518590750013
$ws.Range("J9").Value = 0.003538518590750013
$ws.Range("M9").Value = 4.426331333333334
$ws.Range("N9").Value = 13.278994
$ws.Range("O9").Value = 0.06262029324154843
$ws.Range("P9").Value = 0.06262029324154841
$ws.Range("Q9").Value = 1.414597951826
$ws.Range("R9").Value = 12.731381566434
$ws.Range("S9").Value = 0.0002215830717934365
$ws.Range("T9").Value = 0.0002215830717934365
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.115957
$ws.Range("H10").Value = 0.347871
$ws.Range("I10").Value = 0.001283894527085267
$ws.Range("J10").Value = 0.001283894527085267
$ws.Range("M10").Value = 41.60833666666667
$ws.Range("N10").Value = 124.82501
$ws.Range("O10").Value = 0.5886423873735626
$ws.Range("P10").Value = 0.5886423873735626
$ws.Range("Q10").Value = 4.824777894856666
$ws.Range("R10").Value = 43.42300105371
$ws.Range("S10").Value = 0.0007557547395593226
$ws.Range("T10").Value = 0.0007557547395593227
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.115957
$ws.Range("H11").Value = 0.347871
$ws.Range("I11").Value = 0.001283894527085267
$ws.Range("J11").Value = 0.001283894527085267
$ws.Range("O11").Value = 0.2506218293658061
$ws.Range("P11").Value = 0.2506218293658061
$ws.Range("Q11").Value = 2.054209292823666
$ws.Range("R11").Value = 18.487883635413
$ws.Range("S11").Value = 0.0003217719950908561
$ws.Range("T11").Value = 0.0003217719950908562
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.115957
$ws.Range("H12").Value = 0.347871
$ws.Range("I12").Value = 0.001283894527085267
$ws.Range("J12").Value = 0.001283894527085267
$ws.Range("M12").Value = 6.935318333333334
$ws.Range("N12").Value = 20.805955
$ws.Range("O12").Value = 0.09811549001908282
$ws.Range("P12").Value = 0.09811549001908279
$ws.Range("Q12").Value = 0.8041987079783334
$ws.Range("R12").Value = 7.237788371805
$ws.Range("S12").Value = 0.0001259699406577896
$ws.Range("T12").Value = 0.0001259699406577896
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.115957
$ws.Range("H13").Value = 0.347871
$ws.Range("I13").Value = 0.001283894527085267
$ws.Range("J13").Value = 0.001283894527085267
$ws.Range("M13").Value = 4.426331333333334
$ws.Range("N13").Value = 13.278994
$ws.Range("O13").Value = 0.06262029324154843
$ws.Range("P13").Value = 0.06262029324154841
$ws.Range("Q13").Value = 0.5132641024193333
$ws.Range("R13").Value = 4.619376921774
$ws.Range("S13").Value = 0.00008039785177729855
$ws.Range("T13").Value = 0.00008039785177729855
$ws.Range("G14").Value = 89.26880233333334
$ws.Range("H14").Value = 267.806407
$ws.Range("I14").Value = 0.9883985163053822
$ws.Range("J14").Value = 0.9883985163053823
$ws.Range("M14").Value = 41.60833666666667
$ws.Range("N14").Value = 124.82501
$ws.Range("O14").Value = 0.5886423873735626
$ws.Range("P14").Value = 0.5886423873735626
$ws.Range("Q14").Value = 3714.326381315453
$ws.Range("R14").Value = 33428.93743183908
$ws.Range("S14").Value = 0.5818132623144873
$ws.Range("T14").Value = 0.5818132623144874
$ws.Range("G15").Value = 89.26880233333334
$ws.Range("H15").Value = 267.806407
$ws.Range("I15").Value = 0.9883985163053822
$ws.Range("J15").Value = 0.9883985163053823
$ws.Range("O15").Value = 0.2506218293658061
$ws.Range("P15").Value = 0.2506218293658061
$ws.Range("Q15").Value = 1581.420727617758
$ws.Range("R15").Value = 14232.78654855982
$ws.Range("S15").Value = 0.2477142442989035
$ws.Range("T15").Value = 0.2477142442989035
$ws.Range("G16").Value = 89.26880233333334
$ws.Range("H16").Value = 267.806407
$ws.Range("I16").Value = 0.9883985163053822
$ws.Range("J16").Value = 0.9883985163053823
$ws.Range("M16").Value = 6.935318333333334
$ws.Range("N16").Value = 20.805955
$ws.Range("O16").Value = 0.09811549001908282
$ws.Range("P16").Value = 0.09811549001908279
$ws.Range("Q16").Value = 619.1075614170762
$ws.Range("R16").Value = 5571.968052753686
$ws.Range("S16").Value = 0.09697720476143699
$ws.Range("T16").Value = 0.09697720476143698
$ws.Range("G17").Value = 89.26880233333334
$ws.Range("H17").Value = 267.806407
$ws.Range("I17").Value = 0.9883985163053822
$ws.Range("J17").Value = 0.9883985163053823
$ws.Range("M17").Value = 4.426331333333334
$ws.Range("N17").Value = 13.278994
$ws.Range("O17").Value = 0.06262029324154843
$ws.Range("P17").Value = 0.06262029324154841
$ws.Range("Q17").Value = 395.1332968571731
$ws.Range("R17").Value = 3556.199671714559
$ws.Range("S17").Value = 0.06189380493055442
$ws.Range("T17").Value = 0.06189380493055441
